# Weekly fruit/vegetable price update:
#  - Insert two new observation rows (new rows 335-336) ahead of the
#    existing row that used to be row 335, shifting the remaining
#    price history down by two rows.
#  - Populate the two new rows with the latest "Zapallo italiano"
#    (Primera/Segunda) price readings.
#  - Correct the unit-of-sale / Kg-o-Unidades values for the former
#    row 351 (now row 353), which had been mis-tagged as "70 unidades"
#    even though its quality grade is "Segunda" (paired with "100
#    unidades" everywhere else in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank rows before the current row 335 -------------
$ws.Rows("335:336").Insert()

# --- 2. Fill in the new row 335 (Primera) -----------------------------
$ws.Range("A335").Value2 = 1
$ws.Range("B335").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C335").Value2 = "Arica y Parinacota"
$ws.Range("D335").Value2 = 44900
$ws.Range("E335").Value2 = 15
$ws.Range("F335").Value2 = 100112032
$ws.Range("G335").Value2 = "Zapallo italiano"
$ws.Range("H335").Value2 = "Huracán"
$ws.Range("I335").Value2 = "Primera"
$ws.Range("J335").Value2 = 250
$ws.Range("K335").Value2 = 5000
$ws.Range("L335").Value2 = 6000
$ws.Range("M335").Value2 = 5400
$ws.Range("N335").Value2 = "$/caja 70 unidades"
$ws.Range("O335").Value2 = "Región de Arica y Parinacota"
$ws.Range("P335").Value2 = 77
$ws.Range("Q335").Value2 = 70
$ws.Range("R335").Value2 = "Hortaliza"

# --- 3. Fill in the new row 336 (Segunda) -----------------------------
$ws.Range("A336").Value2 = 1
$ws.Range("B336").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C336").Value2 = "Arica y Parinacota"
$ws.Range("D336").Value2 = 44900
$ws.Range("E336").Value2 = 15
$ws.Range("F336").Value2 = 100112032
$ws.Range("G336").Value2 = "Zapallo italiano"
$ws.Range("H336").Value2 = "Huracán"
$ws.Range("I336").Value2 = "Segunda"
$ws.Range("J336").Value2 = 160
$ws.Range("K336").Value2 = 3500
$ws.Range("L336").Value2 = 4000
$ws.Range("M336").Value2 = 3812
$ws.Range("N336").Value2 = "$/caja 100 unidades"
$ws.Range("O336").Value2 = "Región de Arica y Parinacota"
$ws.Range("P336").Value2 = 38
$ws.Range("Q336").Value2 = 100
$ws.Range("R336").Value2 = "Hortaliza"

# --- 4. Correct the unit-of-sale tagging on the row that shifted ------
#        from 351 to 353 (Segunda -> "100 unidades", not "70 unidades")
$ws.Range("N353").Value2 = "$/caja 100 unidades"
$ws.Range("Q353").Value2 = 100
